$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift H1/I1 (order_by / rem) left into F1/G1, then clear the now-vacated
# H1/I1 cells. This removes the two "is_locked_lbl" / "is_enabled_lbl"
# template columns (and their now-unreferenced shared strings).
$h1 = $ws.Cells.Item(1, 8).Value2
$i1 = $ws.Cells.Item(1, 9).Value2

$ws.Cells.Item(1, 6).Value2 = $h1
$ws.Cells.Item(1, 7).Value2 = $i1

$ws.Cells.Item(1, 8).Clear()
$ws.Cells.Item(1, 9).Clear()
